# Fill NA values in "Current Outside Temp (C)" (column L) and
# "Cups of Coffee Each Day" (column M) with the column mean values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tempValue = 3.638888888888889
$coffeeValue = 0.5888888888888889

$lRows = @(4,7,16,18,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,45,47,53,55,57,58,59,60,61,62,64,65,66,67,68,69,70,71,73,74,75,76,77,78,79,80,81,82,83,84,85,93,94,95,96,97,98,99,101,102,103)
$mRows = @(4,8,38,47,53,60,64,73,74,78,89,102,103)

foreach ($r in $lRows) {
    $ws.Range("L$r").Value = $tempValue
}

foreach ($r in $mRows) {
    $ws.Range("M$r").Value = $coffeeValue
}
